$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column V (phi) to hold the material id ("mat"),
# shifting the existing phi/c columns one to the right (V->W, W->X).
$ws.Columns("V").Insert()

# New header for the inserted column.
$ws.Range("V1").Value = "mat"

# Populate material ids per slice row, grouped by the material's phi/c values.
$ws.Range("V2:V12").Value = 1
$ws.Range("V13:V15").Value = 2
$ws.Range("V16:V21").Value = 3
